$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying "mutable" data (Fecha, Volumen, Unidad de comercialización,
# Precio $/Kg, Kg / unidad) for rows 2-10 has been reshuffled across rows,
# while the rest of each row (Mercado, Región, Producto, Categoría, Variedad,
# Calidad, Precio mínimo/máximo/promedio, Origen, etc.) stays put.
# Apply the final values directly per row.

$rows = @{
    2  = @{ D = 44327; M = 60;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 }
    3  = @{ D = 44306; M = 80;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 }
    4  = @{ D = 44323; M = 80;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 }
    5  = @{ D = 44330; M = 60;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 }
    6  = @{ D = 44309; M = 80;  Q = "`$/caja 14 kilos granel";    S = 821;   T = 14 }
    7  = @{ D = 44322; M = 60;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 }
    8  = @{ D = 44316; M = 120; Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 }
    9  = @{ D = 44302; M = 80;  Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 }
    10 = @{ D = 44313; M = 120; Q = "`$/caja 10 kilos empedrada"; S = 11500; T = 1 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 13).Value = $vals.M  # M: Volumen
    $ws.Cells.Item($r, 17).Value = $vals.Q  # Q: Unidad de comercialización
    $ws.Cells.Item($r, 19).Value = $vals.S  # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $vals.T  # T: Kg / unidad
}
